$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the existing (quote-prefixed text) formats of row 2 down to the new row 3
# before writing values, so no new cell styles get minted.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats

# New row: DealId 337 / IdNumber 0312046845086
$ws.Range("B3").Value = "'337"

# Existing row 2's IdNumber was replaced with a new value
$ws.Range("A2").Value = "'6309108015081"

$ws.Range("A3").Value = "'0312046845086"

# Leave the selection on A2, matching the saved workbook state
$ws.Range("A2").Select() | Out-Null
